$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1. Insert a new column before the table (B) and re-anchor the table to C2:F4
$ws.Columns("B:B").Insert()
$lo.Resize($ws.Range("C2:F4"))
$lo.ListColumns.Item(4).CalculatedColumnFormula = "=(D3-E3)/1000"

# 2. Column widths
$ws.Columns("B:B").ColumnWidth = 23.21875

# 3. New "Min Support" header + values
$ws.Range("B2").Value2 = "Min Support"
$ws.Range("B3").Value2 = 3

# 4. Re-worked Eclat benchmark data
$ws.Range("C3").Value2 = "x1"
$ws.Range("D3").Value2 = 4251
$ws.Range("E3").Value2 = 4159
$ws.Range("C4").Value2 = "x10"
$ws.Range("D4").Value2 = 50631
$ws.Range("E4").Value2 = 38178
$ws.Range("F3:F4").Formula = "=(D3-E3)/1000"

# 5. Merge the Min Support cells (must happen after the table no longer owns column B)
$ws.Range("B3:B4").Merge()
Write-Host "merged B3:B4"

# 6. Styling: "Min Support" header (bold white on black, white bottom border)
$hdr = $ws.Range("B2")
$hdr.Font.Name = "Calibri"
$hdr.Font.Size = 16
$hdr.Font.Bold = $true
$hdr.Font.ThemeColor = 2
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.Interior.Pattern = 1
$hdr.Interior.ThemeColor = 1
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = -4138
$hdr.Borders.Item(9).Color = 16777215

# 7. Styling: Min Support value cell (white Arial on blue)
$val = $ws.Range("B3:B4")
$val.Font.Name = "Arial"
$val.Font.Size = 14
$val.Font.ThemeColor = 2
$val.HorizontalAlignment = -4108
$val.VerticalAlignment = -4108
$val.Interior.Pattern = 1
$val.Interior.ThemeColor = 5

# 8. Styling: "Grandezza Dataset" data cells (x1/x10) -> Arial 17
$ds = $ws.Range("C3:C4")
$ds.Font.Name = "Arial"
$ds.Font.Size = 17

# 9. Styling: numeric data cells -> Arial 12
$num = $ws.Range("D3:F4")
$num.Font.Name = "Arial"
$num.Font.Size = 12
$num.Font.ThemeColor = 1

# 10. Selection / active cell to match saved state
$ws.Range("D4").Select()
